$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching the bold style used by the existing headers (A1:C1)
$ws.Range("D1").Value = "Cellular Module Code Version"
$ws.Range("E1").Value = "Date"
$ws.Range("D1:E1").Font.Bold = $true

# Set column D width to match the diff (target stored width 25.6640625).
# Note: the engine quantizes ColumnWidth to steps of 1/6 (offset 5/6), so the
# closest reproducible stored width is 25.666666... ; 24.8 lands in that bucket.
$ws.Columns.Item(4).ColumnWidth = 24.8

# Add data row
$ws.Range("D2").Value = "2_6"

# Update selection to match target (activeCell D3, sqref D3)
$ws.Range("D3").Select()
